# Adding the changes we made on may 9th
#
# The "falling" sample window stored in Sheet1 is replaced with a newer
# capture: the sensor readings (ax,ay,az,gx,gy,gz) in rows 2-21 are
# overwritten with a fresh 20-sample window, and 10 additional samples
# are appended as new rows 22-31 (continuing the timestamp series and
# the "falling" label). The timestamp (A) / label (B) columns for the
# pre-existing rows 2-21 are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: overwrite the ax,ay,az,gx,gy,gz columns (C:H) for the
#     existing rows 2-21 with the new data window ---
$chData = @(
    @(-2.872348141670227, 9.113740730285643, -1.157746517658234, -0.007177666760981, 0.0068722339347004, 0.0154243474826216),
    @(-2.973343849182129, 9.098371124267578, -1.108406949043274, -0.0125227374956011, 0.0532979927957057, -0.0287106670439243),
    @(-3.070270323753357, 9.036239910125733, -1.239402884244919, -0.0233655963093042, 0.0858265683054924, 0.0064140851609408),
    @(-3.125558137893677, 8.97771692276001, -1.272288262844086, -0.0406225398182868, -0.0181732401251792, 0.0048869219608604),
    @(-3.148761796951294, 8.941205215454101, -1.218037897348404, 0.0109955742955207, -0.015118914656341, 0.0439822971820831),
    @(-3.052603721618652, 8.98315830230713, -1.222838401794433, 0.0197004042565822, -0.0032070425804704, 0.0224492978304624),
    @(-3.015902495384216, 9.072346258163453, -1.212291812896728, 0.0274889357388019, 0.0001527163112768, 0.0108428578823804),
    @(-3.051667261123657, 9.041837882995607, -1.251752722263336, 0.0339030213654041, 0.0074830991216003, -0.0473420582711696),
    @(-3.03490149974823, 8.900234460830688, -1.314584940671921, 0.0368046313524246, 0.0103847095742821, -0.0549778714776039),
    @(-2.981666660308838, 8.602373123168945, -1.842280006408692, 0.0694859251379966, -0.2252565622329712, 0.1149953827261924),
    @(-3.369193434715273, 7.965732479095457, -2.964559459686284, -0.1847867369651794, -0.4732678532600403, 0.2860376536846161),
    @(-3.835544395446778, 7.139311599731442, -4.718908405303957, 0.2762638032436371, -0.4453207552433014, 0.3353650271892547),
    @(-3.949358296394353, 6.115968346595764, -5.519280028343203, 0.3854559659957886, -0.3795000314712524, 0.1736384481191635),
    @(-6.036023616790772, 5.773025989532471, -6.258669376373291, 0.2231185287237167, -0.2501493096351623, -0.5458080768585205),
    @(-5.327067399024963, 5.096813631057738, -6.893664145469667, 1.827708840370178, 2.010510206222534, -3.499036073684692),
    @(-5.896905612945559, 2.338003730773914, -5.030807352066033, -0.6444628238677979, 1.021672129631042, 0.2196060568094253),
    @(-1.743859815597524, 8.946796488761919, 5.909892201423672, -2.739730596542358, 0.935998260974884, -3.772856473922729),
    @(-1.42204113006592, 18.078337097168, 5.143125724792469, -0.8069530129432678, 1.669800162315369, -0.4091269969940185),
    @(1.667438387870793, 5.428582191467266, 1.3269031047821, 0.0166460778564214, -0.098960168659687, 0.0592539273202419),
    @(0.7794665813445962, 4.995001316070575, -3.324089527130142, -0.0222965814173221, 0.0245873257517814, -0.1383609771728515)
)

$startRow = 2
for ($i = 0; $i -lt $chData.Count; $i++) {
    $rowVals = $chData[$i]
    $arr = New-Object "object[,]" 1,6
    for ($j = 0; $j -lt 6; $j++) {
        $arr[0,$j] = $rowVals[$j]
    }
    $r = $startRow + $i
    $ws.Range("C$r`:H$r").Value = $arr
}

# --- Step 2: append 10 new rows (22-31) with full timestamp/label/sensor data ---
$newRows = @(
    @(2000, "falling", 3.068378543853783, 7.068177723884584, -4.896264028549195, 0.0693332031369209, 0.1020144969224929, 0.0342084541916847),
    @(2100, "falling", 3.911060714721674, 7.258105087280274, -4.906297397613526, 0.050854530185461, 0.3074179291725158, 0.2310597747564315),
    @(2200, "falling", 3.277458310127258, 7.035428524017331, -4.093589890003198, 0.0531452745199203, -0.1438587605953216, -0.0389426611363887),
    @(2300, "falling", 3.418862581253053, 6.974825143814088, -4.04072964191437, 0.0024434609804302, 0.0180205255746841, 0.06643159687519069),
    @(2400, "falling", 3.381656622886657, 7.120937466621399, -4.581331634521484, -0.0386372283101081, 0.0261144898831844, 0.2403754740953445),
    @(2500, "falling", 2.801538753509521, 7.392663955688477, -4.599455451965333, 0.0062613687478005, 0.1797471046447754, 0.1189660057425499),
    @(2600, "falling", 2.732593083381657, 7.733071899414061, -4.365746712684627, 0.0025961773935705, -0.09178250283002851, -0.0105374250560998),
    @(2700, "falling", 3.0636076927185, 7.479286289215088, -3.930359029769903, 0.0099265603348612, 0.0328340083360672, 0.058643065392971),
    @(2800, "falling", 2.46818006038666, 7.359129667282104, -4.541467189788817, -0.0207694191485643, -0.0707076489925384, -0.09605856239795681),
    @(2900, "falling", 3.081884574890136, 7.212435054779053, -4.413347625732421, 0.0001527163112768, -0.0478002056479454, -0.06643159687519069)
)

$startRow2 = 22
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowVals = $newRows[$i]
    $arr = New-Object "object[,]" 1,8
    for ($j = 0; $j -lt 8; $j++) {
        $arr[0,$j] = $rowVals[$j]
    }
    $r = $startRow2 + $i
    $ws.Range("A$r`:H$r").Value = $arr
}
